$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("sex") gets a new "U" value for rows 2-6
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 4).Value = "U"
}

# Clear the now-unused sire/dam columns (F:G) for rows 2-6
$ws.Range("F2:G6").ClearContents()

# Move the active selection to D7
$ws.Range("D7").Select()
